$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Row 3 (VT187-0443 / "Call UUID method"): link target + validated page title text
$ws.Range("G3").Value = "wait(2);`nvalidate1;`nlink_Click(generic_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0443_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;"
$ws.Range("H3").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=PB and RE2.2 Semi Auto Frame Work : Generic`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0443`n};`nvalidate4`n{`nvalidate_SystemProperties=uuid`n};"

# Row 4 (VT187-0444 / "OEMInfo Method"): link target + validated page title text
$ws.Range("G4").Value = "wait(2);`nvalidate1;`nlink_Click(generic_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0444_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;"
$ws.Range("H4").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=PB and RE2.2 Semi Auto Frame Work : Generic`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0444`n};`nvalidate4`n{`nvalidate_SystemProperties=oeminfo`n};"

# Update the active selection on the TestCases sheet to A2
$ws.Activate()
$ws.Range("A2").Select()

# Row heights grow to fit the new (longer) wrapped text
$ws.Rows.Item(3).RowHeight = 217.5
$ws.Rows.Item(4).RowHeight = 217.5
